# Update column C (Fitness) values on Sheet1 rows 2-189 to reflect the
# new run's fitness progression, per the commit diff. Rows 190-252 are
# left untouched because their values already equal the new tail value
# (7293) and the diff shows no change for them.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ranges = @(
    @(2, 5, 10768),
    @(6, 6, 10462),
    @(7, 10, 10411),
    @(11, 11, 10062),
    @(12, 15, 9208),
    @(16, 16, 8962),
    @(17, 39, 8140),
    @(40, 57, 8112),
    @(58, 58, 8036),
    @(59, 62, 7748),
    @(63, 68, 7581),
    @(69, 189, 7293)
)

foreach ($r in $ranges) {
    $startRow = $r[0]
    $endRow = $r[1]
    $value = $r[2]
    $rangeAddress = "C" + $startRow + ":C" + $endRow
    $ws.Range($rangeAddress).Value = $value
}
